$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell that used to hold "fff" at C3 moves down to C4 and its text
# is extended with "changechangechange".
$ws.Range("C3").ClearContents()
$ws.Range("C4").Value = "fffchangechangechange"

# The cell that used to hold "ff" at D4 moves down to D6, text unchanged.
$ws.Range("D4").ClearContents()
$ws.Range("D6").Value = "ff"

# Update the active selection to E4.
$ws.Range("E4").Select()
